$d = $word.ActiveDocument

# Update the date heading (first paragraph, unique text, scope to doc start)
$dateRange = $d.Range(0, $d.Paragraphs.Item(1).Range.End)
$dateRange.Find.Execute("2024-01-16 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-17 Wednesday", 2) | Out-Null

# Update each table cell (row-major, 20 rows x 5 cols)
# Scope each Find/Replace to a Document.Range(start,end) built from the cell bounds
# so the replacement only touches that one cell, even when two cells share identical text.
$t = $d.Tables.Item(1)

$cell = $t.Cell(1,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("19+28=47", $true, $false, $false, $false, $false, $true, 1, $false, "91-29=62", 2) | Out-Null
$cell = $t.Cell(1,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("65-46=19", $true, $false, $false, $false, $false, $true, 1, $false, "56-49=7", 2) | Out-Null
$cell = $t.Cell(1,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("17+38=55", $true, $false, $false, $false, $false, $true, 1, $false, "30-22=8", 2) | Out-Null
$cell = $t.Cell(1,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("2+29=31", $true, $false, $false, $false, $false, $true, 1, $false, "92-43=49", 2) | Out-Null
$cell = $t.Cell(1,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("29+14=43", $true, $false, $false, $false, $false, $true, 1, $false, "81-32=49", 2) | Out-Null

$cell = $t.Cell(2,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("59+38=97", $true, $false, $false, $false, $false, $true, 1, $false, "28+47=75", 2) | Out-Null
$cell = $t.Cell(2,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("59+19=78", $true, $false, $false, $false, $false, $true, 1, $false, "70-7=63", 2) | Out-Null
$cell = $t.Cell(2,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("90-55=35", $true, $false, $false, $false, $false, $true, 1, $false, "13+48=61", 2) | Out-Null
$cell = $t.Cell(2,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("25+47=72", $true, $false, $false, $false, $false, $true, 1, $false, "73-24=49", 2) | Out-Null
$cell = $t.Cell(2,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("75-67=8", $true, $false, $false, $false, $false, $true, 1, $false, "49+47=96", 2) | Out-Null

$cell = $t.Cell(3,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("8+56=64", $true, $false, $false, $false, $false, $true, 1, $false, "31-14=17", 2) | Out-Null
$cell = $t.Cell(3,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("90-84=6", $true, $false, $false, $false, $false, $true, 1, $false, "16+8=24", 2) | Out-Null
$cell = $t.Cell(3,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("83-79=4", $true, $false, $false, $false, $false, $true, 1, $false, "37+46=83", 2) | Out-Null
$cell = $t.Cell(3,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("16+7=23", $true, $false, $false, $false, $false, $true, 1, $false, "80-4=76", 2) | Out-Null
$cell = $t.Cell(3,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("60-19=41", $true, $false, $false, $false, $false, $true, 1, $false, "70-22=48", 2) | Out-Null

$cell = $t.Cell(4,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("34-29=5", $true, $false, $false, $false, $false, $true, 1, $false, "52-38=14", 2) | Out-Null
$cell = $t.Cell(4,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("38+53=91", $true, $false, $false, $false, $false, $true, 1, $false, "86-77=9", 2) | Out-Null
$cell = $t.Cell(4,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("9+44=53", $true, $false, $false, $false, $false, $true, 1, $false, "49+23=72", 2) | Out-Null
$cell = $t.Cell(4,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("89+7=96", $true, $false, $false, $false, $false, $true, 1, $false, "6+39=45", 2) | Out-Null
$cell = $t.Cell(4,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("62-35=27", $true, $false, $false, $false, $false, $true, 1, $false, "45-26=19", 2) | Out-Null

$cell = $t.Cell(5,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("86+7=93", $true, $false, $false, $false, $false, $true, 1, $false, "55-28=27", 2) | Out-Null
$cell = $t.Cell(5,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("27+18=45", $true, $false, $false, $false, $false, $true, 1, $false, "73-36=37", 2) | Out-Null
$cell = $t.Cell(5,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("80-34=46", $true, $false, $false, $false, $false, $true, 1, $false, "58-19=39", 2) | Out-Null
$cell = $t.Cell(5,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("70-47=23", $true, $false, $false, $false, $false, $true, 1, $false, "65+18=83", 2) | Out-Null
$cell = $t.Cell(5,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("84-6=78", $true, $false, $false, $false, $false, $true, 1, $false, "73-45=28", 2) | Out-Null

$cell = $t.Cell(6,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("23-8=15", $true, $false, $false, $false, $false, $true, 1, $false, "49+22=71", 2) | Out-Null
$cell = $t.Cell(6,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("75+9=84", $true, $false, $false, $false, $false, $true, 1, $false, "36+45=81", 2) | Out-Null
$cell = $t.Cell(6,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("59+3=62", $true, $false, $false, $false, $false, $true, 1, $false, "82-27=55", 2) | Out-Null
$cell = $t.Cell(6,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("69+9=78", $true, $false, $false, $false, $false, $true, 1, $false, "89+2=91", 2) | Out-Null
$cell = $t.Cell(6,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("28+23=51", $true, $false, $false, $false, $false, $true, 1, $false, "67+16=83", 2) | Out-Null

$cell = $t.Cell(7,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("22-13=9", $true, $false, $false, $false, $false, $true, 1, $false, "15+47=62", 2) | Out-Null
$cell = $t.Cell(7,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("57+7=64", $true, $false, $false, $false, $false, $true, 1, $false, "80-1=79", 2) | Out-Null
$cell = $t.Cell(7,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("61-2=59", $true, $false, $false, $false, $false, $true, 1, $false, "51-18=33", 2) | Out-Null
$cell = $t.Cell(7,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("19+46=65", $true, $false, $false, $false, $false, $true, 1, $false, "39+36=75", 2) | Out-Null
$cell = $t.Cell(7,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("68+8=76", $true, $false, $false, $false, $false, $true, 1, $false, "91-17=74", 2) | Out-Null

$cell = $t.Cell(8,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("70-4=66", $true, $false, $false, $false, $false, $true, 1, $false, "29+16=45", 2) | Out-Null
$cell = $t.Cell(8,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("73-18=55", $true, $false, $false, $false, $false, $true, 1, $false, "72-56=16", 2) | Out-Null
$cell = $t.Cell(8,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("5+76=81", $true, $false, $false, $false, $false, $true, 1, $false, "83-74=9", 2) | Out-Null
$cell = $t.Cell(8,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("83-66=17", $true, $false, $false, $false, $false, $true, 1, $false, "21-7=14", 2) | Out-Null
$cell = $t.Cell(8,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("67+9=76", $true, $false, $false, $false, $false, $true, 1, $false, "63-4=59", 2) | Out-Null

$cell = $t.Cell(9,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("91-56=35", $true, $false, $false, $false, $false, $true, 1, $false, "58+15=73", 2) | Out-Null
$cell = $t.Cell(9,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("67-28=39", $true, $false, $false, $false, $false, $true, 1, $false, "23+59=82", 2) | Out-Null
$cell = $t.Cell(9,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("15+79=94", $true, $false, $false, $false, $false, $true, 1, $false, "61-14=47", 2) | Out-Null
$cell = $t.Cell(9,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("71-47=24", $true, $false, $false, $false, $false, $true, 1, $false, "85-69=16", 2) | Out-Null
$cell = $t.Cell(9,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("43-24=19", $true, $false, $false, $false, $false, $true, 1, $false, "67+8=75", 2) | Out-Null

$cell = $t.Cell(10,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("29+57=86", $true, $false, $false, $false, $false, $true, 1, $false, "6+45=51", 2) | Out-Null
$cell = $t.Cell(10,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("8+86=94", $true, $false, $false, $false, $false, $true, 1, $false, "64-45=19", 2) | Out-Null
$cell = $t.Cell(10,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("7+49=56", $true, $false, $false, $false, $false, $true, 1, $false, "8+24=32", 2) | Out-Null
$cell = $t.Cell(10,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("78+8=86", $true, $false, $false, $false, $false, $true, 1, $false, "48+5=53", 2) | Out-Null
$cell = $t.Cell(10,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("51-39=12", $true, $false, $false, $false, $false, $true, 1, $false, "92-33=59", 2) | Out-Null

$cell = $t.Cell(11,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("39+42=81", $true, $false, $false, $false, $false, $true, 1, $false, "17+34=51", 2) | Out-Null
$cell = $t.Cell(11,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("34+8=42", $true, $false, $false, $false, $false, $true, 1, $false, "4+67=71", 2) | Out-Null
$cell = $t.Cell(11,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("29+56=85", $true, $false, $false, $false, $false, $true, 1, $false, "65-19=46", 2) | Out-Null
$cell = $t.Cell(11,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("26+8=34", $true, $false, $false, $false, $false, $true, 1, $false, "67-29=38", 2) | Out-Null
$cell = $t.Cell(11,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("58+26=84", $true, $false, $false, $false, $false, $true, 1, $false, "29+25=54", 2) | Out-Null

$cell = $t.Cell(12,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("80-13=67", $true, $false, $false, $false, $false, $true, 1, $false, "70-34=36", 2) | Out-Null
$cell = $t.Cell(12,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("68+7=75", $true, $false, $false, $false, $false, $true, 1, $false, "26+55=81", 2) | Out-Null
$cell = $t.Cell(12,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("91-35=56", $true, $false, $false, $false, $false, $true, 1, $false, "25+39=64", 2) | Out-Null
$cell = $t.Cell(12,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("53-44=9", $true, $false, $false, $false, $false, $true, 1, $false, "53-19=34", 2) | Out-Null
$cell = $t.Cell(12,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("56-8=48", $true, $false, $false, $false, $false, $true, 1, $false, "90-31=59", 2) | Out-Null

$cell = $t.Cell(13,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("70-42=28", $true, $false, $false, $false, $false, $true, 1, $false, "16-7=9", 2) | Out-Null
$cell = $t.Cell(13,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("47-8=39", $true, $false, $false, $false, $false, $true, 1, $false, "17+19=36", 2) | Out-Null
$cell = $t.Cell(13,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("90-87=3", $true, $false, $false, $false, $false, $true, 1, $false, "62-14=48", 2) | Out-Null
$cell = $t.Cell(13,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("60-35=25", $true, $false, $false, $false, $false, $true, 1, $false, "27+9=36", 2) | Out-Null
$cell = $t.Cell(13,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("74+9=83", $true, $false, $false, $false, $false, $true, 1, $false, "91-2=89", 2) | Out-Null

$cell = $t.Cell(14,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("6+28=34", $true, $false, $false, $false, $false, $true, 1, $false, "70-29=41", 2) | Out-Null
$cell = $t.Cell(14,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("44+8=52", $true, $false, $false, $false, $false, $true, 1, $false, "47+14=61", 2) | Out-Null
$cell = $t.Cell(14,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("39+3=42", $true, $false, $false, $false, $false, $true, 1, $false, "75+6=81", 2) | Out-Null
$cell = $t.Cell(14,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("79+7=86", $true, $false, $false, $false, $false, $true, 1, $false, "27+7=34", 2) | Out-Null
$cell = $t.Cell(14,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("18+9=27", $true, $false, $false, $false, $false, $true, 1, $false, "38+17=55", 2) | Out-Null

$cell = $t.Cell(15,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("93-4=89", $true, $false, $false, $false, $false, $true, 1, $false, "54-49=5", 2) | Out-Null
$cell = $t.Cell(15,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("9+67=76", $true, $false, $false, $false, $false, $true, 1, $false, "70-27=43", 2) | Out-Null
$cell = $t.Cell(15,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("28+5=33", $true, $false, $false, $false, $false, $true, 1, $false, "72-57=15", 2) | Out-Null
$cell = $t.Cell(15,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("83-18=65", $true, $false, $false, $false, $false, $true, 1, $false, "39+38=77", 2) | Out-Null
$cell = $t.Cell(15,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("37+25=62", $true, $false, $false, $false, $false, $true, 1, $false, "60-25=35", 2) | Out-Null

$cell = $t.Cell(16,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("38+9=47", $true, $false, $false, $false, $false, $true, 1, $false, "91-13=78", 2) | Out-Null
$cell = $t.Cell(16,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("44-16=28", $true, $false, $false, $false, $false, $true, 1, $false, "94-89=5", 2) | Out-Null
$cell = $t.Cell(16,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("4+7=11", $true, $false, $false, $false, $false, $true, 1, $false, "9+84=93", 2) | Out-Null
$cell = $t.Cell(16,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("67+17=84", $true, $false, $false, $false, $false, $true, 1, $false, "72-64=8", 2) | Out-Null
$cell = $t.Cell(16,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("34+27=61", $true, $false, $false, $false, $false, $true, 1, $false, "57+6=63", 2) | Out-Null

$cell = $t.Cell(17,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("41-29=12", $true, $false, $false, $false, $false, $true, 1, $false, "60-32=28", 2) | Out-Null
$cell = $t.Cell(17,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("35+38=73", $true, $false, $false, $false, $false, $true, 1, $false, "84-38=46", 2) | Out-Null
$cell = $t.Cell(17,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("49+38=87", $true, $false, $false, $false, $false, $true, 1, $false, "2+19=21", 2) | Out-Null
$cell = $t.Cell(17,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("70-42=28", $true, $false, $false, $false, $false, $true, 1, $false, "6+69=75", 2) | Out-Null
$cell = $t.Cell(17,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("34-8=26", $true, $false, $false, $false, $false, $true, 1, $false, "38+28=66", 2) | Out-Null

$cell = $t.Cell(18,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("9+56=65", $true, $false, $false, $false, $false, $true, 1, $false, "46+36=82", 2) | Out-Null
$cell = $t.Cell(18,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("46-8=38", $true, $false, $false, $false, $false, $true, 1, $false, "87-79=8", 2) | Out-Null
$cell = $t.Cell(18,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("93-19=74", $true, $false, $false, $false, $false, $true, 1, $false, "59+16=75", 2) | Out-Null
$cell = $t.Cell(18,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("86-57=29", $true, $false, $false, $false, $false, $true, 1, $false, "37+45=82", 2) | Out-Null
$cell = $t.Cell(18,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("18+14=32", $true, $false, $false, $false, $false, $true, 1, $false, "60-45=15", 2) | Out-Null

$cell = $t.Cell(19,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("85-16=69", $true, $false, $false, $false, $false, $true, 1, $false, "73-29=44", 2) | Out-Null
$cell = $t.Cell(19,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("64-28=36", $true, $false, $false, $false, $false, $true, 1, $false, "90-45=45", 2) | Out-Null
$cell = $t.Cell(19,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("48-29=19", $true, $false, $false, $false, $false, $true, 1, $false, "3+79=82", 2) | Out-Null
$cell = $t.Cell(19,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("8+25=33", $true, $false, $false, $false, $false, $true, 1, $false, "48+15=63", 2) | Out-Null
$cell = $t.Cell(19,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("8+78=86", $true, $false, $false, $false, $false, $true, 1, $false, "27+48=75", 2) | Out-Null

$cell = $t.Cell(20,1)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("42-8=34", $true, $false, $false, $false, $false, $true, 1, $false, "52-25=27", 2) | Out-Null
$cell = $t.Cell(20,2)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("28+3=31", $true, $false, $false, $false, $false, $true, 1, $false, "39+9=48", 2) | Out-Null
$cell = $t.Cell(20,3)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("94-88=6", $true, $false, $false, $false, $false, $true, 1, $false, "50-35=15", 2) | Out-Null
$cell = $t.Cell(20,4)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("65-57=8", $true, $false, $false, $false, $false, $true, 1, $false, "38+39=77", 2) | Out-Null
$cell = $t.Cell(20,5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("35+48=83", $true, $false, $false, $false, $false, $true, 1, $false, "77-38=39", 2) | Out-Null
